# Weekly price-data refresh: a new observation row is inserted at row 65
# of the "Hortaliza, Feria Lagunitas de Puerto Montt - Ajo" sheet, pushing
# all subsequent rows (65-188) down by one (to 66-189).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 65; Excel shifts rows 65..188 down to
# 66..189 and extends the used range to A1:R189 automatically.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new price observation.
$ws.Cells.Item(65, 1).Value = 4
$ws.Cells.Item(65, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(65, 3).Value = 'Los Lagos'
$ws.Cells.Item(65, 4).Value = 44540
$ws.Cells.Item(65, 5).Value = 10
$ws.Cells.Item(65, 6).Value = 100112003
$ws.Cells.Item(65, 7).Value = 'Ajo'
$ws.Cells.Item(65, 8).Value = 'Chino'
$ws.Cells.Item(65, 9).Value = 'Primera'
$ws.Cells.Item(65, 10).Value = 240
$ws.Cells.Item(65, 11).Value = 21000
$ws.Cells.Item(65, 12).Value = 22000
$ws.Cells.Item(65, 13).Value = 21500
$ws.Cells.Item(65, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(65, 15).Value = 'China'
$ws.Cells.Item(65, 16).Value = 2150
$ws.Cells.Item(65, 17).Value = 10
$ws.Cells.Item(65, 18).Value = 'Hortaliza'
